$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.000000000000001376127323043572
$ws.Range("E2").Value = 0.000000000000001376127323043572

$ws.Range("D3").Value = 0.9999999999999969
$ws.Range("E3").Value = 0.9999999999999969

$ws.Range("D5").Value = 0.5255307017843521
$ws.Range("E5").Value = 0.4744692982156479

$ws.Range("D7").Value = 0.000000009730464143880366
$ws.Range("E7").Value = 0.9999999902695359

$ws.Range("D8").Value = 0.0001501716521676274
$ws.Range("E8").Value = 0.9998498283478324
$ws.Range("F8").Value = 215.6007537841797
